# Planning.xlsx edit script
# Implements: "[Build 404] [Update] Continuité de l'intro, loin d'être fini
#  Ajout des 4 flèches + de la particule d'explode"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------------
# Row 19 : "Team Blui" / "Nouveau profil" step — progress 20% -> 70%,
#          comment replaced, last-touched date set (now a real date instead
#          of a free-text comment), row no longer needs the taller 30pt row.
# ---------------------------------------------------------------------------
$ws.Range("E19").Value2 = 0.7
$ws.Range("F19").Value2 = "Faire d'autres cubes qui viennent s'écraser à tour de role"

$ws.Range("G19").NumberFormat = "d-mmm"
$ws.Range("G19").Value2 = 41189

$ws.Rows(19).AutoFit()

# ---------------------------------------------------------------------------
# Row 21 : progress 30% -> 70%, comment replaced, last-touched date set as a
#          real date value, row no longer needs the taller 60pt row.
# ---------------------------------------------------------------------------
$ws.Range("E21").Value2 = 0.7
$ws.Range("F21").Value2 = "Ok jusqu'à l'explosion de la particule question 1"

$ws.Range("G21").NumberFormat = "d-mmm"
$ws.Range("G21").Value2 = 41187

$ws.Rows(21).AutoFit()

# ---------------------------------------------------------------------------
# Row 22 : comment replaced, date pushed back a few days
# ---------------------------------------------------------------------------
$ws.Range("C22").Value2 = "Entrée dans le jeu et trans"
$ws.Range("G22").Value2 = 41192

# ---------------------------------------------------------------------------
# Rows 23-25 : these three sub-tasks are now finished -> apply the usual
#              "done" green fill (as used on row 4 etc.), set progress to
#              100% and clear the (no longer meaningful) date.
# ---------------------------------------------------------------------------
$ws.Range("B4:D4").Copy()
$ws.Range("B23:D23").PasteSpecial(-4122)
$ws.Range("B23:D23").PasteSpecial(-4122)
$ws.Range("B24:D24").PasteSpecial(-4122)
$ws.Range("B25:D25").PasteSpecial(-4122)

$ws.Range("E20").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("E25").PasteSpecial(-4122)

$ws.Range("E23").Value2 = 1
$ws.Range("E24").Value2 = 1
$ws.Range("E25").Value2 = 1

$ws.Range("G23").ClearContents()
$ws.Range("G24").ClearContents()
$ws.Range("G25").ClearContents()

# ---------------------------------------------------------------------------
# Row 27 : comment replaced, row now needs the taller 30pt row to fit it.
# ---------------------------------------------------------------------------
$ws.Range("F27").Value2 = "Mettre la config du globaloffset second + les touches dans un txt + petit debugging général"
$ws.Rows(27).RowHeight = 30

# ---------------------------------------------------------------------------
# View state : scrolled up a bit, selection moved to F23
# ---------------------------------------------------------------------------
$ws.Range("F23").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
